# Apply the localization + new "Weight/Single GPU All Layers" column edit
# to the DeepSeek v3 prefill performance workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename worksheet (tab name) to English
# ---------------------------------------------------------------------
$ws.Name = "Performance Analysis"

# ---------------------------------------------------------------------
# 2. Translate title + header row (row 1, row 3) to English
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Performance Analysis Report: deepseek_v3 (EXTEND)"

$ws.Range("A3").Value = "Operator Name"
$ws.Range("B3").Value = "Type"
$ws.Range("H3").Value = "Input"
$ws.Range("I3").Value = "Output"
$ws.Range("J3").Value = "Weight"
$ws.Range("K3").Value = "Compute(us)"
$ws.Range("L3").Value = "Memory(us)"
$ws.Range("M3").Value = "Transfer(us)"
$ws.Range("N3").Value = "Single Layer Latency(us)"
$ws.Range("O3").Value = "Total Time(ms)"
$ws.Range("P3").Value = "Percent(%)"

# ---------------------------------------------------------------------
# 3. Add the new column Q: header + per-row weight bytes,
#    reusing formatting from existing neighbour columns.
# ---------------------------------------------------------------------

# Header cell Q3 -> copy the header style from P3, then set its text.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = "Weight/Single GPU All Layers"

# Data cells Q4:Q20 -> copy the plain numeric style used by columns C:G.
$ws.Range("C4:C20").Copy()
$ws.Range("Q4:Q20").PasteSpecial(-4122)

$ws.Range("Q4").Value = 923467776
$ws.Range("Q5").Value = 575668224
$ws.Range("Q6").Value = 255852544
$ws.Range("Q7").Value = 1790967808
$ws.Range("Q8").Value = 792723456
$ws.Range("Q9").Value = 396361728
$ws.Range("Q10").Value = 425721856
$ws.Range("Q11").Value = 1702887424
$ws.Range("Q12").Value = 851443712
$ws.Range("Q13").Value = 1702887424
$ws.Range("Q14").Value = 851443712
$ws.Range("Q15").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("Q20").Value = 0

# New column needs an explicit width (12 chars), matching the diff's <cols> addition.
$ws.Columns.Item(17).ColumnWidth = 11.1

# ---------------------------------------------------------------------
# 4. Translate the summary labels (rows 24-33) to English
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "Compute Time (ms)"
$ws.Range("A25").Value = "Memory Time (ms)"
$ws.Range("A26").Value = "Transfer Time (ms)"
$ws.Range("A27").Value = "Total Time (ms)"

$ws.Range("A30").Value = "Performance Bottleneck"
$ws.Range("B30").Value = "combine (Total Time: 80.136 ms)"

$ws.Range("A33").Value = "Throughput TPS"

# ---------------------------------------------------------------------
# 5. Append two new summary rows (34, 35) with GPU memory metrics.
# ---------------------------------------------------------------------

# Row 34 + 35 label/value cells reuse the formatting of row 33 (label s=7, value s=8).
$ws.Range("A33:B33").Copy()
$ws.Range("A34:B34").PasteSpecial(-4122)
$ws.Range("A35:B35").PasteSpecial(-4122)

$ws.Range("A34").Value = "Weight Memory/Single GPU (GB)"
$ws.Range("B34").Value = 9.564

$ws.Range("A35").Value = "KV Cache Memory/Single GPU (GB)"
$ws.Range("B35").Value = 0.000033
$ws.Range("B35").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 6. Extend the title merge from A1:O1 to A1:Q1 to span the new column.
# ---------------------------------------------------------------------
$ws.Range("A1:O1").UnMerge()
$ws.Range("A1:Q1").Merge()
